# Commit message: "Added tab for non-standard names"
#
# The "Small N w zero BV" sheet holds paired y/dose values under the
# standard column headers "y" / "dose". This adds a sibling sheet,
# "Small N+zero+names", placed immediately after it, that carries the
# same data but under non-standard column headers "yVar" / "doseVar" -
# useful for exercising the app's handling of custom column names.

$wb = $excel.ActiveWorkbook

# Source data lives on "Small N w zero BV" (6th tab).
$src = $wb.Worksheets.Item("Small N w zero BV")

# Insert the new sheet right after the source sheet.
$newWs = $wb.Worksheets.Add($null, $src)
$newWs.Name = "Small N+zero+names"

# Non-standard header names replacing "y" / "dose".
$newWs.Range("A1").Value = "yVar"
$newWs.Range("B1").Value = "doseVar"

# Copy the data rows verbatim from the source sheet.
$lastRow = $src.Cells.Item($src.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $newWs.Cells.Item($r, 1).Value = $src.Cells.Item($r, 1).Value()
    $newWs.Cells.Item($r, 2).Value = $src.Cells.Item($r, 2).Value()
}

# Match the saved selection/active-cell state of the new tab.
[void]$newWs.Range("B2").Select()
